# "Add files via upload" -- fill in the previously-blank "Asistido Copilot"
# (column D) hour values for the Fase 4 and Fase 5 blocks of the estimate
# sheet, then leave the view scrolled/selected where the author ended up.
# The dependent SUM() totals (D29, D34, D35, D36, D37) are plain formulas,
# so Excel recalculates them automatically after the writes below.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Fase 4: Garantia de Calidad y Depuracion (rows 25-28)
$ws.Range("D25").Value = 1
$ws.Range("D26").Value = 4
$ws.Range("D27").Value = 3
$ws.Range("D28").Value = 4

# Fase 5: Documentacion del proyecto (rows 32-33)
$ws.Range("D32").Value = 2
$ws.Range("D33").Value = 3

# Match the author's final viewport/selection (topLeftCell A20, cell D34
# selected) when the workbook was saved.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 20
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D34").Select()
